$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = "K3mjs"
$ws.Range("H3").Value = "kiz1I"
$ws.Range("H4").Value = "NTOGu"
$ws.Range("H5").Value = "DJIUF"
$ws.Range("H6").Value = "I0oVm"
$ws.Range("H7").Value = "K1rfs"
$ws.Range("H8").Value = "HJZxQ"
$ws.Range("H9").Value = "kEzuT"
$ws.Range("H10").Value = "Alz89"
$ws.Range("H11").Value = "OjxK9"
$ws.Range("H12").Value = "tL0WP"
$ws.Range("H13").Value = "_li9K"
$ws.Range("H14").Value = "aNCVm"
